$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each written cell to remain plain Text (matches source workbook,
# which stores these as inlineStr) by prefixing with a literal quote and
# then resetting the style Excel tags on for the quote-prefix so no stray
# number-format override is left behind on the cell.
$ws.Range("D2").Value = "'29.696.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.67%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.607.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.82%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.70%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'212.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.08%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +1.44%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.67%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'28.08"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.50%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.34%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.18%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.05%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.11%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.596.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.62%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +3.83%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'29.730.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.81%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.51%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'64.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'241.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.11%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.97%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0699"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.57%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.32%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.40%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.02%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'155.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.29%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'15.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.46%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +1.05%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.59%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.93%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.32%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.46%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.71%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.428.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.03%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +3.85%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.25%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -1.05%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.20%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.80%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'56.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.00%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.548"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.48%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +5.98%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.818"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.95%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.95%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.60%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'66.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.81%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.981"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +17.49%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'5.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.23%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.746.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.96%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'86.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.27%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +4.97%  "
$ws.Range("E51").Style = "Normal"
